$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 1.353609635284556
$ws.Range("C2").Value2 = 0.2539931740461725
$ws.Range("D2").Value2 = 0.0262438978229369
$ws.Range("E2").Value2 = 0.0859856760044746
$ws.Range("F2").Value2 = 0.6471769034444748
$ws.Range("I2").Value2 = 0.5584063081530779
$ws.Range("L2").Value2 = 0.2280261146695466
$ws.Range("O2").Value2 = 2.192062247822236

$ws.Range("B3").Value2 = 1.211335594582863
$ws.Range("C3").Value2 = 0.234463823412284
$ws.Range("D3").Value2 = 0.02439403277239194
$ws.Range("E3").Value2 = 0.08708918564981616
$ws.Range("F3").Value2 = 0.647734377875544
$ws.Range("I3").Value2 = 0.5688561124922522
$ws.Range("L3").Value2 = 0.2173551475487585
$ws.Range("O3").Value2 = 2.208821801407851

$ws.Range("B4").Value2 = 1.123850794926682
$ws.Range("C4").Value2 = 0.2224143719674885
$ws.Range("D4").Value2 = 0.02325291009245234
$ws.Range("E4").Value2 = 0.08782035368917374
$ws.Range("F4").Value2 = 0.6486532643167564
$ws.Range("I4").Value2 = 0.575753142022279
$ws.Range("L4").Value2 = 0.2109065978476252
$ws.Range("O4").Value2 = 2.22106646405939

$ws.Range("B5").Value2 = 1.088170069108514
$ws.Range("C5").Value2 = 0.2174897871223038
$ws.Range("D5").Value2 = 0.02278659348616685
$ws.Range("E5").Value2 = 0.08813178773652641
$ws.Range("F5").Value2 = 0.649172405027997
$ws.Range("I5").Value2 = 0.5786843933572214
$ws.Range("L5").Value2 = 0.2083048654136803
$ws.Range("O5").Value2 = 2.226546576890826

$ws.Range("B6").Value2 = 1.082243568008209
$ws.Range("C6").Value2 = 0.2166712077235786
$ws.Range("D6").Value2 = 0.02270908442063302
$ws.Range("E6").Value2 = 0.08818431501452118
$ws.Range("F6").Value2 = 0.6492673378069185
$ws.Range("I6").Value2 = 0.5791784054670046
$ws.Range("L6").Value2 = 0.2078744289042334
$ws.Range("O6").Value2 = 2.227486127074968

$ws.Range("B7").Value2 = 1.123369710529289
$ws.Range("C7").Value2 = 0.2223480148875012
$ws.Range("D7").Value2 = 0.02324662640304354
$ws.Range("E7").Value2 = 0.08782449923345137
$ws.Range("F7").Value2 = 0.6486596802007867
$ws.Range("I7").Value2 = 0.5757921857016122
$ws.Range("L7").Value2 = 0.2108714041644646
$ws.Range("O7").Value2 = 2.221138386974729

$ws.Range("B8").Value2 = 1.30458138523602
$ws.Range("C8").Value2 = 0.2472717602742307
$ws.Range("D8").Value2 = 0.02560718275297802
$ws.Range("E8").Value2 = 0.08635503779971909
$ws.Range("F8").Value2 = 0.6472492646653762
$ws.Range("I8").Value2 = 0.5619094616711102
$ws.Range("L8").Value2 = 0.2243253383029042
$ws.Range("O8").Value2 = 2.197434682087447

$ws.Range("B9").Value2 = 1.658842196621094
$ws.Range("C9").Value2 = 0.2956717304811605
$ws.Range("D9").Value2 = 0.03019299205465131
$ws.Range("E9").Value2 = 0.08389892303558888
$ws.Range("F9").Value2 = 0.6490731178696763
$ws.Range("I9").Value2 = 0.5385112720372263
$ws.Range("L9").Value2 = 0.2515272424181489
$ws.Range("O9").Value2 = 2.166508631736832

$ws.Range("B10").Value2 = 1.918370087783671
$ws.Range("C10").Value2 = 0.3309281768760854
$ws.Range("D10").Value2 = 0.0335345537012941
$ws.Range("E10").Value2 = 0.08235398010184447
$ws.Range("F10").Value2 = 0.653232828331177
$ws.Range("I10").Value2 = 0.5236673728869725
$ws.Range("L10").Value2 = 0.2720110635472963
$ws.Range("O10").Value2 = 2.153343904008125

$ws.Range("B11").Value2 = 2.036258689377291
$ws.Range("C11").Value2 = 0.346898922706032
$ws.Range("D11").Value2 = 0.03504847107448938
$ws.Range("E11").Value2 = 0.08170754492397592
$ws.Range("F11").Value2 = 0.6557422041478276
$ws.Range("I11").Value2 = 0.5174271531601526
$ws.Range("L11").Value2 = 0.2814379857903049
$ws.Range("O11").Value2 = 2.149445806838742

$ws.Range("B12").Value2 = 2.080873502717282
$ws.Range("C12").Value2 = 0.3529366265018723
$ws.Range("D12").Value2 = 0.03562083679940287
$ws.Range("E12").Value2 = 0.0814708664268089
$ws.Range("F12").Value2 = 0.6567815295646255
$ws.Range("I12").Value2 = 0.5151380844476385
$ws.Range("L12").Value2 = 0.2850233025405799
$ws.Range("O12").Value2 = 2.148271556569227

$ws.Range("B13").Value2 = 2.071266143728508
$ws.Range("C13").Value2 = 0.3516367524358657
$ws.Range("D13").Value2 = 0.03549760913796263
$ws.Range("E13").Value2 = 0.08152147846121238
$ws.Range("F13").Value2 = 0.6565537250154563
$ws.Range("I13").Value2 = 0.5156277819944535
$ws.Range("L13").Value2 = 0.284250450436474
$ws.Range("O13").Value2 = 2.148511006604338

$ws.Range("B14").Value2 = 2.039929732063626
$ws.Range("C14").Value2 = 0.3473958520712586
$ws.Range("D14").Value2 = 0.03509557861637091
$ws.Range("E14").Value2 = 0.08168791061572378
$ws.Range("F14").Value2 = 0.6558259228511076
$ws.Range("I14").Value2 = 0.5172373460159214
$ws.Range("L14").Value2 = 0.2817326409034564
$ws.Range("O14").Value2 = 2.149343143416644

$ws.Range("B15").Value2 = 2.020731685121461
$ws.Range("C15").Value2 = 0.3447968562198866
$ws.Range("D15").Value2 = 0.03484920225250931
$ws.Range("E15").Value2 = 0.08179091180542386
$ws.Range("F15").Value2 = 0.6553917338313582
$ws.Range("I15").Value2 = 0.5182328918214942
$ws.Range("L15").Value2 = 0.2801924317597582
$ws.Range("O15").Value2 = 2.149892201169649

$ws.Range("B16").Value2 = 1.910662131980985
$ws.Range("C16").Value2 = 0.3298830608542858
$ws.Range("D16").Value2 = 0.03343548846573441
$ws.Range("E16").Value2 = 0.08239736100688333
$ws.Range("F16").Value2 = 0.6530812815730513
$ws.Range("I16").Value2 = 0.5240855206254906
$ws.Range("L16").Value2 = 0.2713971722761102
$ws.Range("O16").Value2 = 2.153640828946038

$ws.Range("B17").Value2 = 1.843092332736603
$ws.Range("C17").Value2 = 0.320716367833171
$ws.Range("D17").Value2 = 0.0325666148807926
$ws.Range("E17").Value2 = 0.08278383987983418
$ws.Range("F17").Value2 = 0.6518221918253104
$ws.Range("I17").Value2 = 0.5278073509043821
$ws.Range("L17").Value2 = 0.2660293499349393
$ws.Range("O17").Value2 = 2.156476834326668

$ws.Range("B18").Value2 = 1.804211924772517
$ws.Range("C18").Value2 = 0.3154375875065227
$ws.Range("D18").Value2 = 0.03206628248371857
$ws.Range("E18").Value2 = 0.08301143743210559
$ws.Range("F18").Value2 = 0.6511560611011902
$ws.Range("I18").Value2 = 0.5299962560730229
$ws.Range("L18").Value2 = 0.2629521599107676
$ws.Range("O18").Value2 = 2.158304732819488

$ws.Range("B19").Value2 = 1.79104500166568
$ws.Range("C19").Value2 = 0.3136492049085859
$ws.Range("D19").Value2 = 0.03189678007412056
$ws.Range("E19").Value2 = 0.08308940900076855
$ws.Range("F19").Value2 = 0.6509404828863268
$ws.Range("I19").Value2 = 0.5307456523999576
$ws.Range("L19").Value2 = 0.261912037600041
$ws.Range("O19").Value2 = 2.158957373690669

$ws.Range("B20").Value2 = 1.850286930839275
$ws.Range("C20").Value2 = 0.3216928367819776
$ws.Range("D20").Value2 = 0.0326591681581192
$ws.Range("E20").Value2 = 0.08274214943283376
$ws.Range("F20").Value2 = 0.6519502121555192
$ws.Range("I20").Value2 = 0.5274061643684114
$ws.Range("L20").Value2 = 0.2665997046959347
$ws.Range("O20").Value2 = 2.156154569603189

$ws.Range("B21").Value2 = 2.049134748026574
$ws.Range("C21").Value2 = 0.3486417834199926
$ws.Range("D21").Value2 = 0.03521369001565233
$ws.Range("E21").Value2 = 0.08163880523213329
$ws.Range("F21").Value2 = 0.6560372756996316
$ws.Range("I21").Value2 = 0.5167625682729877
$ws.Range("L21").Value2 = 0.2824717613930687
$ws.Range("O21").Value2 = 2.149090521954292

$ws.Range("B22").Value2 = 2.178934613860577
$ws.Range("C22").Value2 = 0.366195615560116
$ws.Range("D22").Value2 = 0.0368778285526048
$ws.Range("E22").Value2 = 0.08096499492719644
$ws.Range("F22").Value2 = 0.6592277865956504
$ws.Range("I22").Value2 = 0.5102376349358835
$ws.Range("L22").Value2 = 0.2929356403276842
$ws.Range("O22").Value2 = 2.14623388967749

$ws.Range("B23").Value2 = 2.109673261190039
$ws.Range("C23").Value2 = 0.3568323061825538
$ws.Range("D23").Value2 = 0.03599015070271605
$ws.Range("E23").Value2 = 0.08132029094270621
$ws.Range("F23").Value2 = 0.6574773157259415
$ws.Range("I23").Value2 = 0.5136805567428446
$ws.Range("L23").Value2 = 0.2873426131604475
$ws.Range("O23").Value2 = 2.147597064933734

$ws.Range("B24").Value2 = 1.84703435761935
$ws.Range("C24").Value2 = 0.3212514023267659
$ws.Range("D24").Value2 = 0.03261732732742217
$ws.Range("E24").Value2 = 0.08276098084618333
$ws.Range("F24").Value2 = 0.6518921543430238
$ws.Range("I24").Value2 = 0.5275873876430559
$ws.Range("L24").Value2 = 0.2663418197601715
$ws.Range("O24").Value2 = 2.15629965055615

$ws.Range("B25").Value2 = 1.563130837478639
$ws.Range("C25").Value2 = 0.2826305307322912
$ws.Range("D25").Value2 = 0.02895717090262195
$ws.Range("E25").Value2 = 0.08451779714566143
$ws.Range("F25").Value2 = 0.6480861064782246
$ws.Range("I25").Value2 = 0.5444300298351301
$ws.Range("L25").Value2 = 0.2440808119951896
$ws.Range("O25").Value2 = 2.173201939449029

